$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column E is treated as text so GPA values like "2.6995" stay as strings
$ws.Range("E2:E18").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'Semester 1'
$ws.Cells.Item(2, 2).Value = 'Database Development'
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 'B'
$ws.Cells.Item(2, 5).Value = '2.6995'

$ws.Cells.Item(3, 1).Value = 'Semester 1'
$ws.Cells.Item(3, 2).Value = 'English for Tertiary Studies'
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 'B+'
$ws.Cells.Item(3, 5).Value = '2.6995'

$ws.Cells.Item(4, 1).Value = 'Semester 1'
$ws.Cells.Item(4, 2).Value = 'Problem Solving and Programming'
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 'B+'
$ws.Cells.Item(4, 5).Value = '2.6995'

$ws.Cells.Item(5, 1).Value = 'Semester 1'
$ws.Cells.Item(5, 2).Value = 'Introduction to Cybersecurity'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 'B-'
$ws.Cells.Item(5, 5).Value = '2.6995'

$ws.Cells.Item(6, 1).Value = 'Semester 1'
$ws.Cells.Item(6, 2).Value = 'System Analysis and Design'
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 4).Value = 'B'
$ws.Cells.Item(6, 5).Value = '2.6995'

$ws.Cells.Item(7, 1).Value = 'Semester 1'
$ws.Cells.Item(7, 2).Value = 'Integrity and Anti-Corruption'
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 'A'
$ws.Cells.Item(7, 5).Value = '2.6995'

$ws.Cells.Item(8, 1).Value = 'Semester 1'
$ws.Cells.Item(8, 2).Value = 'Calculus and Algebra'
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 'F'
$ws.Cells.Item(8, 5).Value = '2.6995'

$ws.Cells.Item(9, 1).Value = 'Semester 2'
$ws.Cells.Item(9, 2).Value = 'Fundamentals of Computer Network'
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 4).Value = 'B+'
$ws.Cells.Item(9, 5).Value = '2.1000'

$ws.Cells.Item(10, 1).Value = 'Semester 2'
$ws.Cells.Item(10, 2).Value = 'Probability and Statistics'
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 'C+'
$ws.Cells.Item(10, 5).Value = '2.1000'

$ws.Cells.Item(11, 1).Value = 'Semester 2'
$ws.Cells.Item(11, 2).Value = 'Computer Architecture'
$ws.Cells.Item(11, 3).Value = 3
$ws.Cells.Item(11, 4).Value = 'C+'
$ws.Cells.Item(11, 5).Value = '2.1000'

$ws.Cells.Item(12, 1).Value = 'Semester 2'
$ws.Cells.Item(12, 2).Value = 'Calculus and Algebra'
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 'F'
$ws.Cells.Item(12, 5).Value = '2.1000'

$ws.Cells.Item(13, 1).Value = 'Semester 3'
$ws.Cells.Item(13, 2).Value = 'Object-Oriented Programming'
$ws.Cells.Item(13, 3).Value = 4
$ws.Cells.Item(13, 4).Value = 'B+'
$ws.Cells.Item(13, 5).Value = '3.4665'

$ws.Cells.Item(14, 1).Value = 'Semester 3'
$ws.Cells.Item(14, 2).Value = 'Introduction to Interface Design'
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 'A-'
$ws.Cells.Item(14, 5).Value = '3.4665'

$ws.Cells.Item(15, 1).Value = 'Semester 3'
$ws.Cells.Item(15, 2).Value = 'Academic English'
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 'A-'
$ws.Cells.Item(15, 5).Value = '3.4665'

$ws.Cells.Item(16, 1).Value = 'Semester 3'
$ws.Cells.Item(16, 2).Value = 'Web-based Integration Systems'
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 'B'
$ws.Cells.Item(16, 5).Value = '3.4665'

$ws.Cells.Item(17, 1).Value = 'Semester 3'
$ws.Cells.Item(17, 2).Value = 'Discrete Mathematics'
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 4).Value = 'A'
$ws.Cells.Item(17, 5).Value = '3.4665'

$ws.Cells.Item(18, 1).Value = 'Semester 3'
$ws.Cells.Item(18, 2).Value = 'Penghayatan Etika dan Peradaban'
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 'B+'
$ws.Cells.Item(18, 5).Value = '3.4665'
